# Nexial commit: "block the display of crypted data variables"
# Adds two new script commands to the #system reference sheet (desktop command list):
#   clickElementOffset(name,xOffset,yOffset)
#   saveTextByLocator(var,locator)
# Both are inserted in alphabetically-sorted order into column G of the hidden
# "#system" sheet, shifting the existing entries down, and the "desktop" named
# range is extended to cover the two additional rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# Re-write column G (rows 33-94) so the two new commands land in their correct
# alphabetically-sorted slots and every other entry shifts down by the
# appropriate amount.
$ws.Range("G33").Value = 'clickElementOffset(name,xOffset,yOffset)'
$ws.Range("G34").Value = 'clickExplorerBar(group,item)'
$ws.Range("G35").Value = 'clickFirstMatchRow(nameValues)'
$ws.Range("G36").Value = 'clickFirstMatchedList(contains)'
$ws.Range("G37").Value = 'clickIcon(label)'
$ws.Range("G38").Value = 'clickList(row)'
$ws.Range("G39").Value = 'clickMenu(menu)'
$ws.Range("G40").Value = 'clickOffset(locator,xOffset,yOffset)'
$ws.Range("G41").Value = 'clickRadio(name)'
$ws.Range("G42").Value = 'clickTab(group,name)'
$ws.Range("G43").Value = 'clickTableCell(row,column)'
$ws.Range("G44").Value = 'clickTableRow(row)'
$ws.Range("G45").Value = 'clickTextPane(name,criteria)'
$ws.Range("G46").Value = 'clickTextPaneRow(var,index)'
$ws.Range("G47").Value = 'closeApplication()'
$ws.Range("G48").Value = 'collapseHierTable()'
$ws.Range("G49").Value = 'editCurrentRow(nameValues)'
$ws.Range("G50").Value = 'editHierCells(var,matchBy,nameValues)'
$ws.Range("G51").Value = 'editTableCells(row,nameValues)'
$ws.Range("G52").Value = 'getRowCount(var)'
$ws.Range("G53").Value = 'hideExplorerBar()'
$ws.Range("G54").Value = 'login(form,username,password)'
$ws.Range("G55").Value = 'maximize()'
$ws.Range("G56").Value = 'minimize()'
$ws.Range("G57").Value = 'resize(width,height)'
$ws.Range("G58").Value = 'saveAllTableRows(var)'
$ws.Range("G59").Value = 'saveAttributeByLocator(var,locator,attribute)'
$ws.Range("G60").Value = 'saveElementCount(var,name)'
$ws.Range("G61").Value = 'saveFirstListData(var,contains)'
$ws.Range("G62").Value = 'saveFirstMatchedListIndex(var,contains)'
$ws.Range("G63").Value = 'saveHierCells(var,matchBy,column,nestedOnly)'
$ws.Range("G64").Value = 'saveHierRow(var,matchBy)'
$ws.Range("G65").Value = 'saveListData(var,contains)'
$ws.Range("G66").Value = 'saveLocatorCount(var,locator)'
$ws.Range("G67").Value = 'saveModalDialogText(var)'
$ws.Range("G68").Value = 'saveModalDialogTextByLocator(var,locater)'
$ws.Range("G69").Value = 'saveProcessId(var,locator)'
$ws.Range("G70").Value = 'saveRowCount(var)'
$ws.Range("G71").Value = 'saveTableRows(var,contains)'
$ws.Range("G72").Value = 'saveTableRowsRange(var,beginRow,endRow)'
$ws.Range("G73").Value = 'saveText(var,name)'
$ws.Range("G74").Value = 'saveTextByLocator(var,locator)'
$ws.Range("G75").Value = 'saveTextPane(var,name,criteria)'
$ws.Range("G76").Value = 'saveWindowTitle(var)'
$ws.Range("G77").Value = 'scanTable(var,name)'
$ws.Range("G78").Value = 'selectCombo(name,text)'
$ws.Range("G79").Value = 'sendKeysToTextBox(name,text1,text2,text3,text4)'
$ws.Range("G80").Value = 'showExplorerBar()'
$ws.Range("G81").Value = 'toggleExplorerBar()'
$ws.Range("G82").Value = 'typeAppendTextArea(name,text1,text2,text3,text4)'
$ws.Range("G83").Value = 'typeAppendTextBox(name,text1,text2,text3,text4)'
$ws.Range("G84").Value = 'typeByLocator(locator,text)'
$ws.Range("G85").Value = 'typeTextArea(name,text1,text2,text3,text4)'
$ws.Range("G86").Value = 'typeTextBox(name,text1,text2,text3,text4)'
$ws.Range("G87").Value = 'useApp(appId)'
$ws.Range("G88").Value = 'useForm(formName)'
$ws.Range("G89").Value = 'useHierTable(var,name)'
$ws.Range("G90").Value = 'useList(var,name)'
$ws.Range("G91").Value = 'useTable(var,name)'
$ws.Range("G92").Value = 'useTableRow(var,row)'
$ws.Range("G93").Value = 'waitFor(name,maxWaitMs)'
$ws.Range("G94").Value = 'waitForLocator(locator,maxWaitMs)'

# Extend the "desktop" named range to include the two newly inserted rows.
$wb.Names.Item("desktop").RefersTo = "='#system'!`$G`$2:`$G`$94"
